$d = $word.ActiveDocument

# Update the date/day heading line.
$d.Content.Find.Execute("2026-02-18 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-02-19 Thursday", 2)

# Update the multiplication-problem grid. The practice grid is one table
# with 5 "problem" rows (1, 5, 10, 15, 20) of 5 columns each, separated by
# blank spacer rows. Because some problem strings repeat verbatim
# elsewhere in the table (e.g. "581×3="), each cell is addressed by its
# explicit (row, column) table coordinates rather than by text search, so
# the correct occurrence is always the one updated.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "606×4="
$t.Cell(1, 2).Range.Text  = "582×7="
$t.Cell(1, 3).Range.Text  = "912×4="
$t.Cell(1, 4).Range.Text  = "921×3="
$t.Cell(1, 5).Range.Text  = "742×6="

$t.Cell(5, 1).Range.Text  = "323×8="
$t.Cell(5, 2).Range.Text  = "559×3="
$t.Cell(5, 3).Range.Text  = "472×8="
$t.Cell(5, 4).Range.Text  = "430×6="
$t.Cell(5, 5).Range.Text  = "554×9="

$t.Cell(10, 1).Range.Text = "775×4="
$t.Cell(10, 2).Range.Text = "621×4="
$t.Cell(10, 3).Range.Text = "522×4="
$t.Cell(10, 4).Range.Text = "486×3="
$t.Cell(10, 5).Range.Text = "669×9="

$t.Cell(15, 1).Range.Text = "935×7="
$t.Cell(15, 2).Range.Text = "531×4="
$t.Cell(15, 3).Range.Text = "292×7="
$t.Cell(15, 4).Range.Text = "527×4="
$t.Cell(15, 5).Range.Text = "589×5="

$t.Cell(20, 1).Range.Text = "658×6="
$t.Cell(20, 2).Range.Text = "124×6="
$t.Cell(20, 3).Range.Text = "104×7="
$t.Cell(20, 4).Range.Text = "678×4="
$t.Cell(20, 5).Range.Text = "847×5="

Write-Host "grid updated"
